$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.593.15'
$ws.Range('E2').Value = '  -4.97%  '
$ws.Range('D3').Value = '1.837.59'
$ws.Range('E3').Value = '  -4.39%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.005'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '314.07'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.50%  '
$ws.Range('E6').Value = '  -0.14%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4212'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -8.40%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3613'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -5.31%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07198'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -7.08%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.8982'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -8.24%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '20.54'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -9.74%  '
$ws.Range('D13').Value = '1.768.67'
$ws.Range('E13').Value = '  -8.52%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.535'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -6.15%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.323'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -6.67%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.06803'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -3.06%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.005'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.21%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '77.47'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -8.35%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000008875'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -6.56%  '
$ws.Range('E20').Value = '  -0.06%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '15.31'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -8.44%  '
$ws.Range('D22').Value = '27.569.85'
$ws.Range('E22').Value = '  -5.08%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.928'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -7.99%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '10.74'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.85%  '
$ws.Range('D25').Value = '1.996.93'
$ws.Range('E25').Value = '  -7.53%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.023'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.33%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '153.24'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.88%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.08'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -4.83%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.288'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -6.14%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '110.53'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -6.28%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.634'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -10.81%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.08864'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -4.75%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.7739'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -10.06%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.489'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -11.92%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.929'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.08%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.003'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.10%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.062'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -14.20%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05314'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -6.61%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.077'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -7.13%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.953'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -5.94%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.01913'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -6.51%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5051'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -8.00%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.733'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -9.26%  '
$ws.Range('B44').Value = 'Algorand'
$ws.Range('C44').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.1628'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -7.18%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.06621'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -4.40%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '8.190'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -12.72%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.4687'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -9.70%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '105.02'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -4.81%  '
$ws.Range('E49').Value = '  -0.13%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '10.14'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -9.75%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.620'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -8.01%  '
